$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data["B"] = @(0.7253364120738013, 0.6826654248632167, 0.6565870575131214, 0.6459915512600674, 0.6442341163874232, 0.6564440333546884, 0.7105987722287352, 0.8177232591540644, 0.8969482396470028, 0.9330938157802677, 0.9467954825485094, 0.9438439726168326, 0.934220783330602, 0.9283281089663831, 0.8945880762683771, 0.8739159690250062, 0.8620359365768593, 0.8580153151664547, 0.8761155212730785, 0.9370469726585213, 0.976951008656556, 0.9556463452742037, 0.8751210891721257, 0.7886486059721847)
$data["D"] = @(0.008979455509845025, 0.007821302310667022, 0.007109052004803829, 0.006818533464379328, 0.006770277192241991, 0.007105135045243571, 0.008580369955929257, 0.01146369786571455, 0.01357566533703647, 0.01453495277009864, 0.0148979864530645, 0.01481981106828556, 0.0145648244362917, 0.01440860757553253, 0.01351294277988302, 0.01296309536711959, 0.01264670161205572, 0.01253955330417966, 0.01302164177634779, 0.01463972651351497, 0.01569590304828239, 0.01513232994977187, 0.0129951738128824, 0.01068476682329589)
$data["E"] = @(0.535411221264134, 0.5303836484837561, 0.5273782259024955, 0.5261743208103127, 0.5259756809625884, 0.5273619048446321, 0.5336609878325973, 0.5466466989470149, 0.5565561050779948, 0.5611406553761427, 0.5628874159754247, 0.5625107497504445, 0.5612841501045551, 0.5605342045994774, 0.5562580084663011, 0.5536541059616624, 0.5521636508488967, 0.5516602623438445, 0.5539305489079851, 0.5616441447109466, 0.5667475832041333, 0.5640182066800108, 0.5538055485846058, 0.5430677134689148)
$data["F"] = @(0.3725575301122603, 0.3579190005688631, 0.3491844527522119, 0.3456888216419358, 0.3451122263728905, 0.3491370512472329, 0.3674575355237906, 0.4053979858665571, 0.4345069338856433, 0.4480188439983834, 0.4531743358452474, 0.4520622819706972, 0.4484422114096702, 0.4462298692768343, 0.4336293326617664, 0.4259684811298996, 0.4215875864890819, 0.4201086579692515, 0.4267813606160047, 0.4495044603318377, 0.4645816106999661, 0.4565139520301926, 0.4264137848138319, 0.3949178138340557)
$data["G"] = @(0.2405468349764988, 0.226671661605323, 0.2183259145934358, 0.2149685448031562, 0.2144136864861821, 0.2182804595433225, 0.2357265785696967, 0.2713219185447855, 0.2983278881917784, 0.3108014760982343, 0.3155521212857622, 0.3145277754261002, 0.3111917696478201, 0.3091519099491364, 0.2975165078587025, 0.290426880171708, 0.2863668463105427, 0.2849952319661924, 0.2911797473744855, 0.3121708975199482, 0.3260482156859013, 0.3186271144183053, 0.2908393266298361, 0.261543210632837)
$data["H"] = @(0.3567772508854716, 0.3535672879066283, 0.3517973117262301, 0.3511266099771149, 0.3510182966112296, 0.3517880615588922, 0.3556287735415395, 0.3647537052407301, 0.3724288176666022, 0.3761312885551433, 0.3775636361793886, 0.3772538074752134, 0.3762485216004023, 0.3756366999462131, 0.3721910951854284, 0.3701313463022728, 0.3689664983619849, 0.3685755140139406, 0.3703485544266272, 0.3765429763793975, 0.3807679937203261, 0.3784968773806696, 0.3702502944206714, 0.3621146240710402)
$data["I"] = @(0.7449121468774171, 0.7561285609061521, 0.7634236498916662, 0.7664991264618539, 0.7670160088778255, 0.7634647110993198, 0.7486948720827185, 0.7229680775906928, 0.706036868273177, 0.6987615982427018, 0.6960679868422757, 0.6966453745764012, 0.6985387631609274, 0.6997065112104028, 0.7065209129375747, 0.7108106462373556, 0.7133181619287976, 0.7141740637101472, 0.7103498388972316, 0.6979809633318599, 0.6902549309903883, 0.6943457266398783, 0.7105580412792705, 0.6997065112104028)
$data["L"] = @(0.1718117218132562, 0.1567895011920371, 0.1475883738917361, 0.1438447805042671, 0.1432235262349622, 0.147537862079048, 0.1666275276720768, 0.2042317180763007, 0.2319524613620416, 0.244581497042148, 0.2493662444995977, 0.2483356611093512, 0.2449750945283853, 0.2429169552229808, 0.2311274764913662, 0.2238996057704981, 0.2197441059010856, 0.2183374395645785, 0.224668842788148, 0.2459621113961674, 0.2598923823176591, 0.2524563587640785, 0.2243210710244199, 0.1940417199848525)
$data["M"] = @(0.1845301833611188, 0.1717276833274681, 0.1638968503600893, 0.1607135292102733, 0.1601854207612803, 0.1638538870039028, 0.1801098250504367, 0.2122148776844881, 0.2359292756922429, 0.2467427694469038, 0.2508410041248439, 0.2499582294799367, 0.2470798672459935, 0.2453172211844077, 0.2352230802143822, 0.2290370223396252, 0.2254814016629823, 0.2242779570305871, 0.2296952887660169, 0.2479252221056498, 0.2598592140246865, 0.2534881151168307, 0.2293976839177674, 0.2035064968160327)
$data["N"] = @(1.572562487964404, 1.558789268756399, 1.550874689875116, 1.547786707418595, 1.547282268483841, 1.550832487221413, 1.567701486266202, 1.605046473129207, 1.635038110963052, 1.649227406568798, 1.65467817712198, 1.653500819708256, 1.649674294891781, 1.647340515732338, 1.634121714676269, 1.626151586075025, 1.621618820441455, 1.620092966998598, 1.626994702512235, 1.650796139409152, 1.666803601923618, 1.658219074196239, 1.626613375727203, 1.594491830780356)
$data["O"] = @(1.1088714286482, 1.072462683276029, 1.050886135341699, 1.0422892150211, 1.040873520315159, 1.050769402002402, 1.096156040069587, 1.191346434927141, 1.265077141312986, 1.29944830712185, 1.312583443202755, 1.309749239858803, 1.300526546915165, 1.294892948334137, 1.262847626534324, 1.243401671755691, 1.232295058410926, 1.228547976226338, 1.245463634415415, 1.303232231101902, 1.341684038602523, 1.32109781722562, 1.244531193322473, 1.164930181577262)

$cols = @("B","D","E","F","G","H","I","L","M","N","O")
$startRow = 2

foreach ($col in $cols) {
    $colValues = $data[$col]
    for ($i = 0; $i -lt $colValues.Length; $i++) {
        $cellRef = "$col$($startRow + $i)"
        $ws.Range($cellRef).Value = $colValues[$i]
    }
}

Write-Host "Updated $($cols.Length) columns x $($colValues.Length) rows"